$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: goto -> ashleyfurniture.com (was newegg.com), waitAfter bumped 2000->7000
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Value = "https://www.ashleyfurniture.com/"
$ws.Cells.Item(2, 8).Value = 7000

# Row 3: click "Sign In / Register" (div) -> click "Login" (span); waitAfter 2000->7000
$ws.Cells.Item(3, 4).Value = "Login"
$ws.Cells.Item(3, 5).Value = "span"
$ws.Cells.Item(3, 8).Value = 7000

# Row 4: waitfortext "Sign In" (div) -> waitfortext "Account Login" (h1)
$ws.Cells.Item(4, 4).Value = "Account Login"
$ws.Cells.Item(4, 5).Value = "h1"

# ---------------------------------------------------------------------------
# New rows 5-13: login flow + post-login assertions
# ---------------------------------------------------------------------------

# Row 5: click email input
$ws.Cells.Item(5, 1).Value = "TC001"
$ws.Cells.Item(5, 2).Value = "Yes"
$ws.Cells.Item(5, 3).Value = "click"
$ws.Cells.Item(5, 4).Value = "email"
$ws.Cells.Item(5, 5).Value = "input"
$ws.Cells.Item(5, 7).Value = 1000
$ws.Cells.Item(5, 8).Value = 2000

# Row 6: type email
$ws.Cells.Item(6, 1).Value = "TC001"
$ws.Cells.Item(6, 2).Value = "Yes"
$ws.Cells.Item(6, 3).Value = "type"
$ws.Cells.Item(6, 4).Value = "email"
$ws.Cells.Item(6, 5).Value = "input"
$ws.Cells.Item(6, 6).Value = "woodsblainem@gmail.com"
$ws.Cells.Item(6, 7).Value = 1000
$ws.Cells.Item(6, 8).Value = 2000

# Row 7: click password input
$ws.Cells.Item(7, 1).Value = "TC001"
$ws.Cells.Item(7, 2).Value = "Yes"
$ws.Cells.Item(7, 3).Value = "click"
$ws.Cells.Item(7, 4).Value = "password"
$ws.Cells.Item(7, 5).Value = "input"
$ws.Cells.Item(7, 7).Value = 1000
$ws.Cells.Item(7, 8).Value = 2000

# Row 8: type password
$ws.Cells.Item(8, 1).Value = "TC001"
$ws.Cells.Item(8, 2).Value = "Yes"
$ws.Cells.Item(8, 3).Value = "type"
$ws.Cells.Item(8, 4).Value = "password"
$ws.Cells.Item(8, 5).Value = "input"
$ws.Cells.Item(8, 6).Value = "Welcome@123456"
$ws.Cells.Item(8, 7).Value = 1000
$ws.Cells.Item(8, 8).Value = 2000

# Row 9: click password input again (before submit)
$ws.Cells.Item(9, 1).Value = "TC001"
$ws.Cells.Item(9, 2).Value = "Yes"
$ws.Cells.Item(9, 3).Value = "click"
$ws.Cells.Item(9, 4).Value = "password"
$ws.Cells.Item(9, 5).Value = "input"
$ws.Cells.Item(9, 6).Value = "Welcome@123456"
$ws.Cells.Item(9, 7).Value = 1000
$ws.Cells.Item(9, 8).Value = 2000

# Row 10: click "Log in" button
$ws.Cells.Item(10, 1).Value = "TC001"
$ws.Cells.Item(10, 2).Value = "Yes"
$ws.Cells.Item(10, 3).Value = "click"
$ws.Cells.Item(10, 4).Value = "Log in"
$ws.Cells.Item(10, 5).Value = "button"
$ws.Cells.Item(10, 7).Value = 1000
$ws.Cells.Item(10, 8).Value = 5000

# Row 11: waitfortext "Welcome" (span)
$ws.Cells.Item(11, 1).Value = "TC001"
$ws.Cells.Item(11, 2).Value = "no"
$ws.Cells.Item(11, 3).Value = "waitfortext"
$ws.Cells.Item(11, 4).Value = "Welcome"
$ws.Cells.Item(11, 5).Value = "span"
$ws.Cells.Item(11, 7).Value = 1000
$ws.Cells.Item(11, 8).Value = 2000

# Row 12: goto cart
$ws.Cells.Item(12, 1).Value = "TC001"
$ws.Cells.Item(12, 2).Value = "no"
$ws.Cells.Item(12, 3).Value = "goto"
$ws.Cells.Item(12, 4).Value = "https://secure.newegg.com/shop/cart"
$ws.Cells.Item(12, 7).Value = 1000
$ws.Cells.Item(12, 8).Value = 5000

# Row 13: waitfortext "Summary" (h3) -- replaces the old blank placeholder row
$ws.Cells.Item(13, 1).Value = "TC001"
$ws.Cells.Item(13, 2).Value = "no"
$ws.Cells.Item(13, 3).Value = "waitfortext"
$ws.Cells.Item(13, 4).Value = "Summary"
$ws.Cells.Item(13, 5).Value = "h3"
$ws.Cells.Item(13, 7).Value = 1000

# ---------------------------------------------------------------------------
# Hyperlinks. Re-point D2 at the new URL, and add new link cells F6/F8/F9.
# Hyperlinks.Add() restyles the target with a brand-new "Hyperlink" cell
# style; copy/paste-special the formats from the existing hyperlink cell
# (D2) afterwards so the cells share the workbook's single Hyperlink xf,
# matching how the sheet already does it.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.ashleyfurniture.com/")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:woodsblainem@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F8"), "mailto:Welcome@123456")
$ws.Hyperlinks.Add($ws.Range("F9"), "mailto:Welcome@123456")

$ws.Range("D2").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Selection moved to B1 (matches the recorded sheetView <selection>)
# ---------------------------------------------------------------------------
$ws.Range("B1").Select()
